$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change, in order, mapped to column letters B, C, D, F, G, K, M, O
$cols = @("B", "C", "D", "F", "G", "K", "M", "O")

# New values for rows 2..25 (one row per array, values ordered per $cols)
$data = @(
    @(7.606131260226796, 5.751856910868507, 5.06192125803937, 25.25945130385455, 3.636275081473431, 7.458445890786765, 19.97654202861501, 22.58566493674111),
    @(7.377235253427808, 5.717459883635561, 5.015845950737418, 25.24974544060774, 3.638115586870246, 7.259909549470789, 19.38498955004004, 22.62847728711807),
    @(7.234382033268284, 5.695938486449563, 4.986785028342326, 25.25022058531678, 3.639304350111979, 7.133610582942855, 19.02068059439236, 22.65951251182187),
    @(7.175682996269472, 5.687071338551586, 4.974753498186038, 25.25203176650866, 3.639803589247196, 7.081079221741129, 18.87219530899819, 22.67335042822848),
    @(7.165909686835303, 5.685593212772677, 4.972744399999689, 25.25243017219965, 3.63988738343262, 7.072293494308892, 18.8475453280499, 22.67572002766371),
    @(7.233592230733925, 5.695819288551258, 4.98662352569182, 25.25023846366398, 3.639311022998422, 7.132906374606243, 19.01867784026386, 22.65969431772965),
    @(7.527742759937049, 5.740082500752775, 5.046197243167567, 25.25476926941091, 3.636897538649521, 7.390927328139709, 19.7729470224061, 22.5994392932262),
    @(8.082378077685485, 5.823526422685578, 5.156704069893439, 25.31468401195566, 3.632627988667851, 7.860343345516509, 21.23312899292488, 22.51908760281704),
    @(8.471658719519946, 5.882562777267091, 5.233779323840584, 25.38970953558893, 3.629770271922773, 8.180927827228578, 22.28169345193771, 22.48326975802066),
    @(8.643949480324915, 5.90888048197591, 5.267896044441742, 25.4305267769183, 3.628530123281358, 8.321133863961434, 22.75109401400452, 22.47204654863183),
    @(8.708440474218602, 5.918765032128113, 5.280674860412804, 25.44693890616376, 3.628069062168646, 8.373390478691297, 22.92757792673921, 22.46852763796697),
    @(8.694585543089572, 5.916639900337841, 5.277929030742651, 25.44336187416424, 3.628167980161467, 8.362173677750533, 22.88962796251472, 22.469252956196),
    @(8.649270593611574, 5.909695343781263, 5.268950205514461, 25.43185790210666, 3.628492020304445, 8.325449953199287, 22.76563980712389, 22.4717423839687),
    @(8.621414237618012, 5.905430884999443, 5.263431993575491, 25.4249356211625, 3.628691617081777, 8.302845941816804, 22.68952344596023, 22.47336248984715),
    @(8.46029715427678, 5.880831679733789, 5.231530306017594, 25.3871761081573, 3.629852518531863, 8.171649394141166, 22.25084876614374, 22.48410542637273),
    @(8.360182549840971, 5.8656001472045, 5.211714302190792, 25.36572076179882, 3.63057998591382, 8.089702719994273, 21.97966009150152, 22.49199603603594),
    @(8.302151249356486, 5.856788947078917, 5.200227956671554, 25.35401014352255, 3.631004041341979, 8.042040711738316, 21.82297390853345, 22.49701167108492),
    @(8.28242790355641, 5.853797088283136, 5.196323782544544, 25.35015348834542, 3.631148588603959, 8.025813277485748, 21.76980709313503, 22.49879176778874),
    @(8.370886761339179, 5.867226817840819, 5.21383296339841, 25.36793956557812, 3.63050196287833, 8.098480979189791, 22.00860294420515, 22.49110666727998),
    @(8.662601559915514, 5.911737368908495, 5.271591351194028, 25.43521101996411, 3.628396610033403, 8.336259500441709, 22.8020938586866, 22.47099132369501),
    @(8.848847716638456, 5.940351355490375, 5.308519168178968, 25.48474263401324, 3.627070490157543, 8.486775969702709, 23.31322627202044, 22.46210671210897),
    @(8.749866732009743, 5.925124411301675, 5.288886645185122, 25.45779975247493, 3.627773720104337, 8.406897691543611, 23.04116146847797, 22.46645807495054),
    @(8.36604885918266, 5.866491569153095, 5.212875408677372, 25.36693449923727, 3.630537218946763, 8.094514038562606, 21.99552028469007, 22.49150725838681),
    @(7.935221158198355, 5.801335220263713, 5.127512921750752, 25.29301828712395, 3.633733759986708, 7.737487795013874, 20.84151004386216, 22.53675875826899)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}
